$d = $word.ActiveDocument

# Change 1: "serve to alleviate " -> "alleviate "
$d.Content.Find.Execute("serve to alleviate ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "alleviate ", 2)

# Change 2: "shortages and counteract inflationary" -> "shortages while counteracting inflationary"
$d.Content.Find.Execute("shortages and counteract inflationary", $true, $false, $false, $false, $false,
                         $true, 1, $false, "shortages while counteracting inflationary", 2)
